$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" footer date from 1/23/14
#    to 9/6/14 everywhere it appears: once on the Slide Master, and
#    once on each of the Custom Layouts (11 of them).
#    NB: use $p.SlideMaster (not $p.Slides.Item(n).Master) to reach
#    the master/layouts - CustomLayouts.Item(N) only resolves the
#    correct layout when walked from Presentation.SlideMaster.
# ------------------------------------------------------------------
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "9/6/14"
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "9/6/14"
        }
    }
}

# ------------------------------------------------------------------
# 2) Slide 1 architecture diagram: rename the "HDFS" source box to
#    "HDFS/S3" and the "ZeroMQ" source box to "Kinesis".
# ------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$group = $slide1.Shapes.Item(1)

for ($i = 1; $i -le $group.GroupItems.Count; $i++) {
    $item = $group.GroupItems.Item($i)
    if ($item.HasTextFrame) {
        $txt = $item.TextFrame.TextRange.Text
        if ($txt -eq "HDFS" -and $item.Name -eq "Rounded Rectangle 47") {
            $item.TextFrame.TextRange.Text = "HDFS/S3"
        } elseif ($txt -eq "ZeroMQ") {
            $item.TextFrame.TextRange.Text = "Kinesis"
        }
    }
}
